$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at row 58 (resistor "1k5 250mW" added for open-drain
#     pull-ups R31 R81 R113 R139). Everything at/after row 58 shifts down. ---
$ws.Rows("58:58").Insert()

# Row 57 keeps its data but grows slightly taller (matches the new row 58).
$ws.Rows("57:57").RowHeight = 15.65
$ws.Rows("58:58").RowHeight = 15.65

# --- Populate the newly inserted row 58 ---
$ws.Range("A58").Value = "1k5 250mW"
$ws.Range("B58").Value = 4
$ws.Range("C58").Value = "R31 R81 R113 R139"
$ws.Range("D58").Value = "Resistor SMD 1206"
$ws.Range("E58").Formula = '=HYPERLINK("http://www.digikey.com/product-search/en?KeyWords="&J58,J58)'
$ws.Range("F58").Formula = '=HYPERLINK("http://ar.mouser.com/Search/Refine.aspx?Keyword="&K58,K58)'
$ws.Range("G58").Formula = "=100*B58"
$ws.Range("H58").Value = 0.1
$ws.Range("I58").Formula = "=B58*H58"
$ws.Range("J58").Value = "CR1206-JW-152ELFCT-ND"
$ws.Range("K58").Value = "652-CR1206JW-152ELF"

# --- Defined names: print area / helper range now cover the extra row ---
$ws.PageSetup.PrintArea = "`$A`$1:`$C`$101"
$wb.Names.Item("ciaa_nxp_1").RefersTo = "='ciaa-nxp'!`$A`$2:`$C`$101"

# --- Reflect where the edit was made in the saved view state ---
$ws.Range("A58").Select()
